$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Short Term"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Short Term")

# Row 120 updates
$ws1.Range("B120").Value = 5.04
$ws1.Range("C120").Value = -4.72
$ws1.Range("D120").Value = 5.87
$ws1.Range("F120").Value = 28.11
$ws1.Range("G120").Value = -17.35

# Row 121 updates
$ws1.Range("B121").Value = -9.46
$ws1.Range("C121").Value = -7.83
$ws1.Range("D121").Value = -1.95

# Row 122 updates
$ws1.Range("B122").Value = 31.39
$ws1.Range("C122").Value = 37.14
$ws1.Range("D122").Value = -0.42

# Row 123 updates
$ws1.Range("B123").Value = 4.95
$ws1.Range("C123").Value = 4.09
$ws1.Range("D123").Value = 21.31

# Row 124 updates
$ws1.Range("B124").Value = 36.78
$ws1.Range("C124").Value = 37.9
$ws1.Range("D124").Value = -7.77

# Row 125 updates
$ws1.Range("B125").Value = -34.25
$ws1.Range("C125").Value = -39.01
$ws1.Range("D125").Value = -12.82
$ws1.Range("E125").Value = 18.34
$ws1.Range("F125").Value = 12.68
$ws1.Range("G125").Value = -1.53

# New row 126
$ws1.Range("A125").Copy()
$ws1.Range("A126").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws1.Range("A126").Value = 45778
$ws1.Range("B126").Value = 5.93
$ws1.Range("C126").Value = 5.02
$ws1.Range("D126").Value = -0.49
$ws1.Range("E126").Value = 34.65
$ws1.Range("F126").Value = 29.34
$ws1.Range("G126").Value = -2.68

# ---------------------------------------------------------------
# Sheet 2: "Medium Term"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Medium Term")

# Row 106 updates
$ws2.Range("B106").Value = 29.33
$ws2.Range("C106").Value = 16.71
$ws2.Range("D106").Value = 12.07

# Row 107 updates
$ws2.Range("B107").Value = 22.26
$ws2.Range("C107").Value = 21.76
$ws2.Range("D107").Value = 14.16

# Row 108 updates
$ws2.Range("B108").Value = 28.43
$ws2.Range("C108").Value = 25.55
$ws2.Range("D108").Value = 15.7

# Row 109 updates
$ws2.Range("C109").Value = 33.78
$ws2.Range("D109").Value = 20.83

# Row 110 updates
$ws2.Range("C110").Value = 45.08
$ws2.Range("D110").Value = 30.06

# Row 111 updates
$ws2.Range("B111").Value = 56.24
$ws2.Range("C111").Value = 43.11
$ws2.Range("D111").Value = 28.56

# New row 112
$ws2.Range("A111").Copy()
$ws2.Range("A112").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("A112").Value = 45778
$ws2.Range("B112").Value = 47.11
$ws2.Range("C112").Value = 43.04
$ws2.Range("D112").Value = 30.24
